# Update scripts/paths for IncrementalProgress off model calculators
# - Renames the existing sheet to "RTP2017"
# - Adds a new "RTP2021" sheet with the new IncrementalProgress (IP) model run rows
# - Updates the absolute path recorded for the workbook
# - Updates selection/view state on both sheets

$wb = $excel.ActiveWorkbook

# --- Rename existing sheet, add the new one right after it -----------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RTP2017"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "RTP2021"

# --- Populate RTP2021 --------------------------------------------------
$ws2.Range("A1").Value = "year"
$ws2.Range("B1").Value = "directory"
$ws2.Range("C1").Value = "category"

$ws2.Range("A2").Value = 2035
$ws2.Range("B2").Value = "2035_TM151_IPA_loPop_loAOC_00"
$ws2.Range("C2").Value = "IP"

$ws2.Range("A3").Value = 2035
$ws2.Range("B3").Value = "2035_TM151_IPA_loPop_hiAOC_00"

$ws2.Range("A4").Value = 2035
$ws2.Range("B4").Value = "2035_TM151_IPA_hiPop_loAOC_00"

$ws2.Range("A5").Value = 2035
$ws2.Range("B5").Value = "2035_TM151_IPA_hiPop_hiAOC_00"

$ws2.Range("C3").Value = "IP_hiAOC"
$ws2.Range("C4").Value = "IP_hiPop"
$ws2.Range("C5").Value = "IP_hiPop_hiAOC"

# --- Formatting: column A centered (matches RTP2017's year column), -------
# --- columns B/C left-aligned -- all using the 10pt data font -------------
$ws2.Range("A1:A5").Font.Size = 10
$ws2.Range("A1:A5").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("B1:C5").Font.Size = 10

# Column widths (fit to content, closest achievable to the authored sizes)
$ws2.Columns.Item(1).ColumnWidth = 7.834333333333333
$ws2.Columns.Item(2).ColumnWidth = 28.66766666666667
$ws2.Columns.Item(3).ColumnWidth = 15.667666666666666

# --- Freeze the header row and set the selection on RTP2021 ---------------
$ws2.Activate()
$ws2.Range("A2").Select()
$ws2.Application.ActiveWindow.FreezePanes = $true
$ws2.Range("C4").Select()

# --- Update the selection on RTP2017 (no longer the active tab) -----------
$ws1.Activate()
$ws1.Range("A10:C13").Select()

# --- RTP2021 is the sheet that was active/selected when the file was saved
$ws2.Activate()
